$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 46, pushing existing rows 46-76 down to 48-78.
$ws.Range("A46:A47").EntireRow.Insert()

# New row 46: Flame Seedless, week of 2022-02-04 (serial 44596)
$ws.Cells.Item(46, 1).Value = 8
$ws.Cells.Item(46, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(46, 3).Value = "Coquimbo"
$ws.Cells.Item(46, 4).Value = 44596
$ws.Cells.Item(46, 5).Value = 4
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100109
$ws.Cells.Item(46, 8).Value = "Uva"
$ws.Cells.Item(46, 9).Value = 100109001
$ws.Cells.Item(46, 10).Value = "Uva"
$ws.Cells.Item(46, 11).Value = "Flame Seedless"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 400
$ws.Cells.Item(46, 14).Value = 8000
$ws.Cells.Item(46, 15).Value = 8500
$ws.Cells.Item(46, 16).Value = 8250
$ws.Cells.Item(46, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(46, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(46, 19).Value = 458
$ws.Cells.Item(46, 20).Value = 18

# New row 47: Superior Seedless, week of 2022-02-04 (serial 44596)
$ws.Cells.Item(47, 1).Value = 8
$ws.Cells.Item(47, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(47, 3).Value = "Coquimbo"
$ws.Cells.Item(47, 4).Value = 44596
$ws.Cells.Item(47, 5).Value = 4
$ws.Cells.Item(47, 6).Value = "Fruta"
$ws.Cells.Item(47, 7).Value = 100109
$ws.Cells.Item(47, 8).Value = "Uva"
$ws.Cells.Item(47, 9).Value = 100109001
$ws.Cells.Item(47, 10).Value = "Uva"
$ws.Cells.Item(47, 11).Value = "Superior Seedless"
$ws.Cells.Item(47, 12).Value = "Primera"
$ws.Cells.Item(47, 13).Value = 300
$ws.Cells.Item(47, 14).Value = 10000
$ws.Cells.Item(47, 15).Value = 11000
$ws.Cells.Item(47, 16).Value = 10500
$ws.Cells.Item(47, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(47, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(47, 19).Value = 583
$ws.Cells.Item(47, 20).Value = 18
